# Update review intersections based on Jingxi reviews.
#
# Column C holds the "Jingxi, Jerome, Ann" intersection titles; rows 4-8
# (papers #2-#6) gained a match in that intersection and are highlighted
# with a red fill. A new row 18 records an additional "Daniel, Jingxi, Ann"
# (column D) intersection title.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 is written first so its shared-string entry is appended before the
# column-C titles below (matches shared string ordering of the source edit).
$ws.Range("D18").Value = "Interface Development for Digital Twin of an Electric Motor Based on Empirical Performance Model"

$ws.Range("C4").Value = "Combining Low-Code Programming and SDL-Based Modeling with Snap! in the Industry 4.0 Context"
$ws.Range("C5").Value = "A Methodology for Digital Twin Modeling and Deployment for Industry 4.0"
$ws.Range("C6").Value = "Digital Twins Driving Model Based on Petri Net in Industrial Pipeline"
$ws.Range("C7").Value = "Automated Model Transformation in modeling Digital Twins of Industrial Internet-of-Things Applications utilizing AutomationML"
$ws.Range("C8").Value = "Proof of Concept for a Roundtrip Engineering IS for the New Enterprise in the Industry 4.0 Era"

# Highlight the newly-matched intersection cells with the reviewer's red fill.
$ws.Range("C4:C8").Interior.Color = 192

# The active cell moved to B4 after the edit.
$ws.Range("B4").Select()
